$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("S3").Value = 2023
